# Apply the style edits described by the commit:
#   "Made grammar use linebreaks within a production, and paragraph-breaks
#    after that. In-progress work on lists..."
#
# Concretely, this touches word/styles.xml:
#   1. "Heading1" (aka "Heading 1#") gains a pageBreakBefore so that each
#      numbered heading/production starts on a new page.
#   2. "Grammar" (aka "Grammar#") shrinks from 9pt to 8pt so more of a
#      grammar production's line-broken body fits together.

$d = $word.ActiveDocument

# 1. Heading1: force a page break before every occurrence of the style.
$heading1 = $d.Styles("Heading1")
$heading1.ParagraphFormat.PageBreakBefore = $true

# 2. Grammar: reduce the font size from 9pt (18 half-points) to 8pt
#    (16 half-points).
$grammar = $d.Styles("Grammar")
$grammar.Font.Size = 8
